# Apply cryptocurrency price/volume updates (GitHub Actions scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "38.012.48"

$ws.Cells.Item(3, 4).Value = "2.056.48"
$ws.Cells.Item(3, 5).Value = "  +1.91%  "

$ws.Cells.Item(4, 5).Value = "  +0.08%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "229.94"
$ws.Cells.Item(5, 5).Value = "  +1.64%  "

$ws.Cells.Item(6, 5).Value = "  +2.61%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "58.15"
$ws.Cells.Item(7, 5).Value = "  +6.37%  "

$ws.Cells.Item(8, 5).Value = "  +0.01%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.387"
$ws.Cells.Item(9, 5).Value = "  +2.96%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.0807"
$ws.Cells.Item(10, 5).Value = "  +2.96%  "

$ws.Cells.Item(12, 4).Value = "2.359.85"

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "14.61"
$ws.Cells.Item(13, 5).Value = "  +3.59%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "20.67"
$ws.Cells.Item(14, 5).Value = "  +2.47%  "

$ws.Cells.Item(15, 5).Value = "  +1.93%  "

$ws.Cells.Item(16, 5).Value = "  +2.82%  "

$ws.Cells.Item(17, 4).Value = "2.053.76"
$ws.Cells.Item(17, 5).Value = "  +1.64%  "

$ws.Cells.Item(18, 4).Value = "37.859.72"
$ws.Cells.Item(18, 5).Value = "  +2.30%  "

$ws.Cells.Item(19, 5).Value = "  -1.09%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "69.77"
$ws.Cells.Item(20, 5).Value = "  +1.24%  "

$ws.Cells.Item(21, 4).Value = "0.0₃0830"
$ws.Cells.Item(21, 5).Value = "  +1.58%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "224.89"
$ws.Cells.Item(22, 5).Value = "  +0.79%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.999"
$ws.Cells.Item(23, 5).Value = "  +0.00%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.46"
$ws.Cells.Item(24, 5).Value = "  +0.66%  "

$ws.Cells.Item(25, 5).Value = "  +3.27%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "9.29"
$ws.Cells.Item(26, 5).Value = "  +1.44%  "

$ws.Cells.Item(27, 5).Value = "  +0.08%  "

$ws.Cells.Item(28, 5).Value = "  +7.96%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "19.05"
$ws.Cells.Item(29, 5).Value = "  +1.91%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.35"
$ws.Cells.Item(30, 5).Value = "  +0.35%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "4.53"
$ws.Cells.Item(32, 5).Value = "  +1.06%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "4.58"
$ws.Cells.Item(33, 5).Value = "  +3.99%  "

$ws.Cells.Item(34, 5).Value = "  +0.09%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "2.02"
$ws.Cells.Item(35, 5).Value = "  +9.40%  "

$ws.Cells.Item(36, 5).Value = "  +0.04%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "5.97"
$ws.Cells.Item(37, 5).Value = "  +12.82%  "

$ws.Cells.Item(38, 5).Value = "  +5.29%  "

$ws.Cells.Item(39, 5).Value = "  -0.07%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "98.04"
$ws.Cells.Item(40, 5).Value = "  +3.32%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.0218"
$ws.Cells.Item(41, 5).Value = "  +1.38%  "

$ws.Cells.Item(42, 4).Value = "1.485.65"
$ws.Cells.Item(42, 5).Value = "  +0.84%  "

$ws.Cells.Item(43, 5).Value = "  +3.28%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.0934"
$ws.Cells.Item(44, 5).Value = "  +2.44%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "16.67"
$ws.Cells.Item(45, 5).Value = "  +2.50%  "

$ws.Cells.Item(46, 2).Value = "TrustWalletToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "1.13"
$ws.Cells.Item(46, 5).Value = "  +0.74%  "

$ws.Cells.Item(47, 2).Value = "FTXToken"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "4.12"
$ws.Cells.Item(47, 5).Value = "  +16.21%  "

$ws.Cells.Item(48, 5).Value = "  +0.93%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "2.97"
$ws.Cells.Item(49, 5).Value = "  +1.79%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "7.05"
$ws.Cells.Item(50, 5).Value = "  -1.86%  "

$ws.Cells.Item(51, 4).Value = "2.248.25"
$ws.Cells.Item(51, 5).Value = "  +2.12%  "
